# Remove the "price weight" column (old column G, "قیمت برای ما وزن")
# from the data sheet, shifting the customer price column (old H) left
# into G, and simplify the F1 header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Delete column G entirely; this shifts column H (and its contents) left
# into the now-vacant G, matching the canonical diff.
$ws.Columns.Item(7).Delete()

# Update the F1 header text.
$ws.Range("F1").Value = "قیمت برای ما"
